$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data dictionary sheet: reorder rows to match the updated data file's
# column order, update a couple of descriptions, and add two new rows
# describing the new "seasonkey" and "known_death" columns.

$ellipsis = [char]0x2026
$seasonkeyDesc = "numeric variable that correponds with seasonid (fall_jun_2009=1, spring_oct_2009=2, " + $ellipsis + ", spring_oct_2020=24)"

$rows = @(
    @("column_name",     "description"),
    @("orderkey",        "unique identifier for each trapping occurrence starting at 1 first day trapping in 2009"),
    @("seasonkey",       $seasonkeyDesc),
    @("seasonid",        "unique identifier for each selected window (season_window_year)"),
    @("season",          "fall or spring trapping season"),
    @("window",          "the select 30-day windows (jun and oct) that will be used for the population demographic model"),
    @("year",            "year of trapping"),
    @("date",             "date in yyyy-mm-dd format "),
    @("julian_date",      "julian date starting each year on 01-01 with 1"),
    @("trap_type",        "in which trap the degu was trapped, Sh=Sherman, T=Tomahawk"),
    @("id",                "unique ear tag number of each degu"),
    @("sex",               "sex (M = male, F = female)"),
    @("radio_collar",      "frequency of radio collar if the degu has one"),
    @("weight",            "weight in g of the degu"),
    @("w_collar",          "weight in g of the degu with collar if it has one"),
    @("status",            "whether trapped for the first time (collection of tissue)"),
    @("age",               "Adult or Pup"),
    @("vulva",             "whether the vulva of the female is closed (0), semi-open (1), open (2)"),
    @("breeding_status",   "breeding status of the females; 0 = nonbreeding, 1= definitely pregnant, 2 = definitely lactating, 3 = definitely postlactating meaning that the female was lactating within this spring season, 4 = both pregnant and/or (post)lactating, 5 = unclear breeding status"),
    @("date_death",        "from the trapping file notes, date of death because degus was found death"),
    @("known_death",       "indicator variable of known death (0 = no known death, 1 = known to have died just before the sampling occasion)"),
    @("notes",             "any notes")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
}

# Scroll / selection state, matching the saved view in the edited workbook.
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
$ws.Range("C23").Select()
